# "Fixed image processing image"
#
# The single slide is a software-architecture diagram that still has the
# generic placeholder captions ("Actor" / "Subsystem") from the template.
# This fills in the real labels, widens the "Line Paramaterization" box so
# its (now longer) caption fits, and shrinks the arrow that follows it so
# it still meets the box's left edge.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Assert-ShapeName($shape, $expected) {
    if ($shape.Name -ne $expected) {
        throw "expected shape '$expected' but found '$($shape.Name)'"
    }
}

# --- Actors (ovals) --------------------------------------------------
$ovalIn = $s.Shapes.Item(1)
Assert-ShapeName $ovalIn "Oval 3"
$ovalIn.TextFrame.TextRange.Text = "Input Image"

$ovalOut = $s.Shapes.Item(2)
Assert-ShapeName $ovalOut "Oval 4"
$ovalOut.TextFrame.TextRange.Text = "Data legible by planner"

# --- Subsystems (rectangles) -----------------------------------------
# "Voxelize Image" ends up as two runs ("Voxelize" + " Image") in the
# real deck. Re-asserting the (already black) font color on the first
# word forces the engine to split the run into two without changing any
# visible formatting, matching that run boundary.
$rectVoxelize = $s.Shapes.Item(3)
Assert-ShapeName $rectVoxelize "Rectangle 5"
$voxelizeRange = $rectVoxelize.TextFrame2.TextRange
$voxelizeRange.Text = "Voxelize Image"
$voxelizeRange.Characters(1, 8).Font.Color.RGB = 0

$rectSegmentation = $s.Shapes.Item(4)
Assert-ShapeName $rectSegmentation "Rectangle 6"
$rectSegmentation.TextFrame.TextRange.Text = "Line-based Segmentation"

# "Line Paramaterization" likewise ends up as two runs ("Line " +
# "Paramaterization"); same trick as above.
$rectParam = $s.Shapes.Item(5)
Assert-ShapeName $rectParam "Rectangle 7"
$paramRange = $rectParam.TextFrame2.TextRange
$paramRange.Text = "Line Paramaterization"
$paramRange.Characters(1, 5).Font.Color.RGB = 0

# The box grows/shifts left to fit the longer caption; only its Left and
# Width change (Top/Height stay as they were).
$rectParam.Left = 654.3529663085938
$rectParam.Width = 97.41448974609375

$rectOccupancy = $s.Shapes.Item(6)
Assert-ShapeName $rectOccupancy "Rectangle 8"
$rectOccupancy.TextFrame.TextRange.Text = "Determine voxel occupancy"

# --- Connector --------------------------------------------------------
# The arrow from "Line-based Segmentation" into "Line Paramaterization"
# shrinks because the latter's left edge moved left.
$paramConnector = $s.Shapes.Item(10)
Assert-ShapeName $paramConnector "Straight Arrow Connector 18"
$paramConnector.Width = 59.92378234863281
